$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("B10").Value = 47438
$ws.Range("C10").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D10").Value = 401.81
$ws.Range("E10").Value = 480.05
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0

# Row 11
$ws.Range("B11").Value = 59408
$ws.Range("C11").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D11").Value = 388.17
$ws.Range("E11").Value = 463.78
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = 2329.02

# Row 27
$ws.Range("F27").Value = 34
$ws.Range("G27").Value = 870.74

# Row 37
$ws.Range("F37").Value = 30
$ws.Range("G37").Value = 768.3

# Row 46
$ws.Range("B46").Value = 24629.95

# Row 55
$ws.Range("F55").Value = 195
$ws.Range("G55").Value = 37613.55

# Row 59
$ws.Range("F59").Value = 87
$ws.Range("G59").Value = 1978.38

# Row 61
$ws.Range("F61").Value = 43
$ws.Range("G61").Value = 1087.47

# Row 66
$ws.Range("F66").Value = 45
$ws.Range("G66").Value = 828.9

# Row 68
$ws.Range("F68").Value = 270
$ws.Range("G68").Value = 25255.8

# Row 84
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0

# Row 85
$ws.Range("B85").Value = 143170.56

# Row 121
$ws.Range("F121").Value = 286
$ws.Range("G121").Value = 32083.48

# Row 123
$ws.Range("F123").Value = 176
$ws.Range("G123").Value = 7902.4

# Row 134
$ws.Range("B134").Value = 78750.37

# Row 136
$ws.Range("F136").Value = 49
$ws.Range("G136").Value = 1764

# Row 138
$ws.Range("F138").Value = 28
$ws.Range("G138").Value = 1475.88

# Row 142
$ws.Range("F142").Value = 89
$ws.Range("G142").Value = 12587.27

# Row 144
$ws.Range("F144").Value = 51
$ws.Range("G144").Value = 4655.79

# Row 147
$ws.Range("F147").Value = 32
$ws.Range("G147").Value = 3239.36

# Row 150
$ws.Range("F150").Value = 11
$ws.Range("G150").Value = 339.46

# Row 151
$ws.Range("F151").Value = 12
$ws.Range("G151").Value = 177.48

# Row 152
$ws.Range("F152").Value = 51
$ws.Range("G152").Value = 2163.93

# Row 156
$ws.Range("F156").Value = 36
$ws.Range("G156").Value = 879.48

# Row 158
$ws.Range("F158").Value = 81
$ws.Range("G158").Value = 1561.68

# Row 159
$ws.Range("B159").Value = 62616.73

# Row 169
$ws.Range("F169").Value = 63
$ws.Range("G169").Value = 3117.24

# Row 170
$ws.Range("F170").Value = 70
$ws.Range("G170").Value = 3463.6

# Row 174
$ws.Range("F174").Value = 215
$ws.Range("G174").Value = 4231.2

# Row 180
$ws.Range("B180").Value = 32818.55

# Row 185
$ws.Range("F185").Value = 18
$ws.Range("G185").Value = 1175.4

# Row 187
$ws.Range("F187").Value = 25
$ws.Range("G187").Value = 2322.25

# Row 198
$ws.Range("B198").Value = 38604.45

# Row 292
$ws.Range("B292").Value = 41864
$ws.Range("C292").Value = "HAM-THERMOSTEEL 1000 ML WITH PLAIN LID"
$ws.Range("F292").Value = 0
$ws.Range("G292").Value = 0

# Row 293
$ws.Range("B293").Value = 56449
$ws.Range("C293").Value = "HAM-Thermosteel 1000 Ml With Plain Lid"
$ws.Range("F293").Value = 24
$ws.Range("G293").Value = 16128.96

# Row 329
$ws.Range("F329").Value = 182
$ws.Range("G329").Value = 6060.6

# Row 343
$ws.Range("F343").Value = 1
$ws.Range("G343").Value = 46.87

# Row 349
$ws.Range("B349").Value = 132669.81

# Row 354
$ws.Range("F354").Value = 11
$ws.Range("G354").Value = 880.22

# Row 355
$ws.Range("F355").Value = 12
$ws.Range("G355").Value = 2073.48

# Row 365
$ws.Range("F365").Value = 18
$ws.Range("G365").Value = 2372.94

# Row 367
$ws.Range("F367").Value = 6
$ws.Range("G367").Value = 794.88

# Row 382
$ws.Range("F382").Value = 6
$ws.Range("G382").Value = 658.92

# Row 399
$ws.Range("F399").Value = 37
$ws.Range("G399").Value = 2171.9

# Row 409
$ws.Range("F409").Value = 111
$ws.Range("G409").Value = 19017.63

# Row 419
$ws.Range("F419").Value = 230
$ws.Range("G419").Value = 9471.4

# Row 421
$ws.Range("F421").Value = 330
$ws.Range("G421").Value = 13117.5

# Row 423
$ws.Range("B423").Value = 107797.5

# Row 425
$ws.Range("F425").Value = 1
$ws.Range("G425").Value = 183.58

# Row 437
$ws.Range("B437").Value = 13803.43

# Row 466
$ws.Range("F466").Value = 64
$ws.Range("G466").Value = 2659.84

# Row 471
$ws.Range("F471").Value = 7
$ws.Range("G471").Value = 294.07

# Row 479
$ws.Range("F479").Value = 16
$ws.Range("G479").Value = 2314.24

# Row 481
$ws.Range("B481").Value = 35785.19

# Row 494
$ws.Range("F494").Value = 51
$ws.Range("G494").Value = 2856.51

# Row 495
$ws.Range("F495").Value = 2
$ws.Range("G495").Value = 331.98

# Row 496
$ws.Range("F496").Value = 84
$ws.Range("G496").Value = 12622.68

# Row 497
$ws.Range("B497").Value = 27119.74

# Row 502
$ws.Range("F502").Value = 12
$ws.Range("G502").Value = 551.64

# Row 511
$ws.Range("B511").Value = 34023.67

# Row 527
$ws.Range("F527").Value = 148
$ws.Range("G527").Value = 5507.08

# Row 532
$ws.Range("B532").Value = 144263.69

# Row 546
$ws.Range("F546").Value = 10
$ws.Range("G546").Value = 322.4

# Row 556
$ws.Range("B556").Value = 13221.99

# Row 559
$ws.Range("F559").Value = 86
$ws.Range("G559").Value = 4347.3

# Row 561
$ws.Range("F561").Value = 774
$ws.Range("G561").Value = 9984.6

# Row 567
$ws.Range("B567").Value = 46699.3

# Row 621
$ws.Range("F621").Value = 185
$ws.Range("G621").Value = 11233.2

# Row 623
$ws.Range("B623").Value = 55667
$ws.Range("C623").Value = "NES-Maggi Atta Noodles Masala 290G"
$ws.Range("D623").Value = 85.76000000000001
$ws.Range("E623").Value = 97.25
$ws.Range("F623").Value = 31
$ws.Range("G623").Value = 2658.56

# Row 624
$ws.Range("B624").Value = 49151
$ws.Range("C624").Value = "NES-MAGGI Atta Noodles Masala 290g"
$ws.Range("D624").Value = 78.09999999999999
$ws.Range("E624").Value = 88.58
$ws.Range("F624").Value = 1
$ws.Range("G624").Value = 78.09999999999999

# Row 638
$ws.Range("B638").Value = 130452.71

# Row 665
$ws.Range("F665").Value = 17
$ws.Range("G665").Value = 3069.86

# Row 667
$ws.Range("B667").Value = 21112.83

# Row 669
$ws.Range("F669").Value = 33
$ws.Range("G669").Value = 2614.59

# Row 671
$ws.Range("F671").Value = 160
$ws.Range("G671").Value = 9904

# Row 688
$ws.Range("B688").Value = 78482.25999999999

# Row 702
$ws.Range("F702").Value = 71
$ws.Range("G702").Value = 3566.33

# Row 711
$ws.Range("B711").Value = 16686.78

# Row 713
$ws.Range("F713").Value = 15
$ws.Range("G713").Value = 1958.25

# Row 720
$ws.Range("B720").Value = 19379.44

# Row 778
$ws.Range("F778").Value = 2
$ws.Range("G778").Value = 143.12

# Row 785
$ws.Range("B785").Value = 12895.76

# Row 831
$ws.Range("F831").Value = 21
$ws.Range("G831").Value = 5742.66

# Row 837
$ws.Range("B837").Value = 178484.61

# Row 843
$ws.Range("F843").Value = 49
$ws.Range("G843").Value = 5331.69

# Row 857
$ws.Range("F857").Value = 352
$ws.Range("G857").Value = 27663.68

# Row 859
$ws.Range("F859").Value = 338
$ws.Range("G859").Value = 34783.58

# Row 861
$ws.Range("F861").Value = 284
$ws.Range("G861").Value = 10459.72

# Row 867
$ws.Range("B867").Value = 184194.08

# Row 923
$ws.Range("B923").Value = 2198802.29

# Row 924
$ws.Range("B924").Value = 2198802.29
